$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price records were inserted into the daily log.
# The first goes right after the existing row 53 (becoming the new row 54),
# pushing the former rows 54-69 down to 55-70.
$ws.Rows(54).Insert()

$ws.Cells.Item(54, 1).Value = 9
$ws.Cells.Item(54, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(54, 3).Value = "Metropolitana"
$ws.Cells.Item(54, 4).Value = 44841
$ws.Cells.Item(54, 5).Value = 13
$ws.Cells.Item(54, 6).Value = 100112035
$ws.Cells.Item(54, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 45
$ws.Cells.Item(54, 11).Value = 17000
$ws.Cells.Item(54, 12).Value = 17000
$ws.Cells.Item(54, 13).Value = 17000
$ws.Cells.Item(54, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(54, 15).Value = "Hijuelas"
$ws.Cells.Item(54, 16).Value = 1133
$ws.Cells.Item(54, 17).Value = 15
$ws.Cells.Item(54, 18).Value = "Hortaliza"

# The second new record is inserted after the (now shifted) row 70,
# becoming the new row 71, pushing the former rows 70-74 down to 72-76.
$ws.Rows(71).Insert()

$ws.Cells.Item(71, 1).Value = 9
$ws.Cells.Item(71, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(71, 3).Value = "Metropolitana"
$ws.Cells.Item(71, 4).Value = 44385
$ws.Cells.Item(71, 5).Value = 13
$ws.Cells.Item(71, 6).Value = 100112035
$ws.Cells.Item(71, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Segunda"
$ws.Cells.Item(71, 10).Value = 16
$ws.Cells.Item(71, 11).Value = 12000
$ws.Cells.Item(71, 12).Value = 12000
$ws.Cells.Item(71, 13).Value = 12000
$ws.Cells.Item(71, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(71, 15).Value = "Hijuelas"
$ws.Cells.Item(71, 16).Value = 800
$ws.Cells.Item(71, 17).Value = 15
$ws.Cells.Item(71, 18).Value = "Hortaliza"
